# Updated cryptos list (price/volume refresh) matching the GitHub Actions
# scraper commit. Cells in column D that look numeric are prefixed with a
# leading apostrophe so Excel keeps them as text (matching how the sheet
# originally stored prices as text, e.g. "64.114.87"), instead of being
# auto-coerced into a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.125.55"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "3.152.31"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("D5").Value = "'592.54"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'146.17"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.143.45"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "'5.90"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "'37.22"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "3.672.34"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "63.946.02"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "3.148.55"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").Value = "'468.42"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").Value = "'14.38"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "'13.02"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("D26").Value = "'81.34"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'9.77"
$ws.Range("E28").Value = "  +8.45%  "
$ws.Range("D29").Value = "'7.42"
$ws.Range("E29").Value = "  +7.94%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'27.74"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").Value = "0.0₃0842"
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'2.32"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'6.17"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  -4.91%  "
$ws.Range("D40").Value = "'463.90"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").Value = "'51.44"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +5.18%  "
$ws.Range("D43").Value = "'0.294"
$ws.Range("E43").Value = "  +5.41%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.929.57"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0373"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "'40.18"
$ws.Range("E46").Value = "  +12.45%  "
$ws.Range("D47").Value = "'0.109"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "'129.14"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("E51").Value = "  -0.75%  "
